# Insert a new parameter row "chemical_recycling_pyrolysis" right after
# "chemical_recycling_gasification" (currently row 9), pushing every row
# below it down by one. Populate the new row's value column with TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 is currently "fossil_routes" - insert a new blank row above it,
# shifting fossil_routes (and everything after) down to row 11.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
